# Reorders the weekly price-record rows (rows 2-24, columns A:R) on the
# active sheet according to a fixed permutation. The header row (row 1) is
# left untouched. Every destination row ends up with the exact same
# A:R content (values + the date column's style) that some particular
# source row had before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 24

# Snapshot every source row's A:R values BEFORE any writes happen, so that
# overlapping reads/writes during the permutation can't clobber data we
# still need to read later.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:R$r").Value()
}

# destinationRow -> sourceRow
$rowMap = @{
    2  = 22
    3  = 10
    4  = 7
    5  = 18
    6  = 23
    7  = 8
    8  = 5
    9  = 3
    10 = 2
    11 = 21
    12 = 4
    13 = 19
    14 = 6
    15 = 16
    16 = 13
    17 = 14
    18 = 20
    19 = 17
    20 = 11
    21 = 24
    22 = 9
    23 = 12
    24 = 15
}

for ($destRow = $firstDataRow; $destRow -le $lastDataRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:R$destRow").Value = $snapshot[$srcRow]
}
